$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1/IF, matching the style already used by the other
# header cells (e.g. H1: bold, bordered, centered) via a format-only copy.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data columns I (I0) and J (IF)
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 2

$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 6

$ws.Range("I4").Value = 9
$ws.Range("J4").Value = 9

$ws.Range("I5").Value = 9
$ws.Range("J5").Value = 9

$ws.Range("I6").Value = 7
$ws.Range("J6").Value = 7

$ws.Range("I7").Value = 4
$ws.Range("J7").Value = 4

$ws.Range("I8").Value = 9
$ws.Range("J8").Value = 9
